$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 74 (shifts existing rows 74..174 down to 75..175),
# copying formatting (incl. the date style) from the row above.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new data point.
$ws.Cells.Item(74, 1).Value = 4
$ws.Cells.Item(74, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(74, 3).Value = "Los Lagos"
$ws.Cells.Item(74, 4).Value = 44467
$ws.Cells.Item(74, 5).Value = 10
$ws.Cells.Item(74, 6).Value = 100112045
$ws.Cells.Item(74, 7).Value = "Zapallo"
$ws.Cells.Item(74, 8).Value = "Paine"
$ws.Cells.Item(74, 9).Value = "1a (guarda)"
$ws.Cells.Item(74, 10).Value = 1100
$ws.Cells.Item(74, 11).Value = 600
$ws.Cells.Item(74, 12).Value = 600
$ws.Cells.Item(74, 13).Value = 600
$ws.Cells.Item(74, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(74, 15).Value = "Región Metropolitana"
$ws.Cells.Item(74, 16).Value = 600
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = "Hortaliza"
